# -----------------------------------------------------------------------
# "Took notes on prior pose estimation"
#
# 1. Prepend a new paragraph containing a pull-quote (3 runs) and move the
#    _GoBack bookmark to sit at the end of that new paragraph.
# 2. Insert a blank paragraph after the quote (before the original
#    "Inverse Depth Maps - ..." paragraph).
# 3. Clean up two stray/empty <w:pict> placeholder runs that were left in
#    the last paragraph (around the two "ρi" occurrences) without
#    otherwise touching the visible text there.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1) Detach the existing _GoBack bookmark -----------------------------
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

# --- 2) Build the new quote paragraph + following blank paragraph --------
$firstPara = $d.Paragraphs(1).Range
$insertPos = $firstPara.Start

$quoteText = "“Empirical studies show that LSD-SLAM, which relies on depth ﬁlters, performs consistently more powerfully than ORB-SLAM” – Mobile SLAM"

# Trailing sentinel "#" lets us add the bookmark exactly at the end of the
# quote text without landing on the (buggy) very-last-character-of-paragraph
# boundary; it is stripped again immediately afterwards.
$firstPara.InsertBefore($quoteText + "#" + "`r" + "`r")

$bookmarkPos = $insertPos + $quoteText.Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$sentinelRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$sentinelRange.Delete()

# --- 3) Remove the two empty <w:pict> runs in the final paragraph --------
# First pict: right after "...sixth parameter, " + italic "ρi"
$anchor1 = $d.Content
$anchor1.Find.Execute("sixth parameter, ", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$afterItalic1 = $d.Range($anchor1.End + 2, $anchor1.End + 2)
$afterItalic1.Find.Execute("ρi , represents", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "ρi , represents", 2)

# Second pict: right after "...depth as 1/" + italic "ρi"
$anchor2 = $d.Content
$anchor2.Find.Execute("depth as 1/", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$afterItalic2 = $d.Range($anchor2.End + 2, $anchor2.End + 2)
$afterItalic2.Find.Execute("1/ρi )", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "1/ρi )", 2)
